$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-0.859180821046573, 1.210703344943181, 1.1198820233157438, 3.0236999925478205)"
$ws.Range("C2").Value = "NIG(1.9418981826723927, 1.6671071541101417, 2.0615484416824525, 5.734867423369177)"
$ws.Range("D2").Value = "EXN(2.7146515958262034, 0.3944902712061403, 1.563268500006644)"
$ws.Range("E2").Value = "NCT(2.176422702269768, 2.183083665956064, -0.009728249792020919, 3.1450405360482776)"
